$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-8 with the new combined tuple-like strings
$ws.Range("A2").Value = "('Beast', ['Token Creature — Beast', '3/3'])"
$ws.Range("A3").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A4").Value = "('Elemental', ['Token Creature — Elemental', 'Flying', '5/5'])"
$ws.Range("A5").Value = "('Goat', ['Token Creature — Goat', '0/1'])"
$ws.Range("A6").Value = "('Goblin Soldier', ['Token Creature — Goblin Soldier', '1/1'])"
$ws.Range("A7").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A8").Value = "('Worm', ['Token Creature — Worm', '1/1'])"

# Remove now-unused rows 9 through 25
$ws.Range("A9:A25").ClearContents()
